$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range
$x1 = '<w:p w14:paraId="34B09502" w14:textId="1FDB8D5D" w:rsidR="009C55E9" w:rsidRPr="005C3DB4" w:rsidRDefault="005C3DB4" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="005C3DB4"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">&lt;draft&gt; AUBRIA.ai Website </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005C3DB4"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Verbage</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p1.InsertXML($x1)

$p5 = $d.Paragraphs(5).Range
$x5 = '<w:p w14:paraId="232381C9" w14:textId="49427316" w:rsidR="005C3DB4" w:rsidRDefault="005C3DB4" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">We are a student-led team that are students in the </w:t></w:r><w:r w:rsidR="00C0690C"><w:t xml:space="preserve">Department of </w:t></w:r><w:r><w:t>Computer Science &amp; Software Engineering</w:t></w:r><w:r w:rsidR="00C0690C"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00AC2ED2"><w:t>(CSSE)</w:t></w:r><w:r><w:t xml:space="preserve"> in the Samuel Ginn College of Engineering at Auburn University. We are passionate about Artificial Intelligence, Agentic AI, and Cybersecurity. We are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>affliated</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> with AI@AU (The Auburn University Artificial Intelligence Initiative, </w:t></w:r><w:hyperlink r:id="rId4" w:history="1"><w:r w:rsidRPr="0014099E"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://eng.auburn.edu/ai-au/</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve">) as well as the Auburn University Center for AI &amp; Cybersecurity Engineering (AU-CAICE, </w:t></w:r><w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidRPr="0014099E"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://eng.auburn.edu/au-caice/</w:t></w:r></w:hyperlink><w:r><w:t>).</w:t></w:r></w:p>'
$p5.InsertXML($x5)

$p7 = $d.Paragraphs(7).Range
$x7 = '<w:p w14:paraId="3F1443ED" w14:textId="447D342C" w:rsidR="005C3DB4" w:rsidRDefault="005C3DB4" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Currently, we are interested in developing AUBRIA Keynote Talks for the Auburn University research community for AU Conferences, Workshop, Meetings, etc.</w:t></w:r><w:r w:rsidR="00E6464A"><w:t xml:space="preserve"> If you are interested in our team developing an AUBRIA Talk for </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00E6464A"><w:t>you</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00E6464A"><w:t xml:space="preserve"> please contact us or our faculty advisor Dr. Gerry Dozier.</w:t></w:r></w:p>'
$p7.InsertXML($x7)

$p13 = $d.Paragraphs(13).Range
$x13 = '<w:p w14:paraId="30C93491" w14:textId="258C0DFD" w:rsidR="009E0DB6" w:rsidRDefault="00793334" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">AUBRIA is the brainchild of Auburn University PSFS Director Christine Cline. </w:t></w:r><w:r w:rsidR="00401272"><w:t xml:space="preserve">During the development of the 2026 AI@AU/Team Sciences: Building Research Communities in AI (Showcase &amp; Workshop), Christine mentioned the idea of having an AI-Generated Keynote Speaker to Dr. Jennifer Kerpelman (Showcase &amp; Workshop </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00401272"><w:t>co-Chair</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00401272"><w:t xml:space="preserve">) who in turn mentioned the idea to Dr. Gerry Dozier (the other Showcase &amp; Workshop </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00401272"><w:t>co-Chair</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00401272"><w:t xml:space="preserve">). </w:t></w:r></w:p>'
$p13.InsertXML($x13)

$p15 = $d.Paragraphs(15).Range
$x15 = '<w:p w14:paraId="2CA20EF2" w14:textId="2C6D9906" w:rsidR="00401272" w:rsidRDefault="00401272" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Dr. Dozier set out to find a team of students from the </w:t></w:r><w:r w:rsidR="00AC2ED2"><w:t xml:space="preserve">Department of Computer Science &amp; Software Engineering (CSSE) student body to develop an AI-Generated Keynote Speaker. Dr. Dozier found us!!! A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00AC2ED2"><w:t>Phd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00AC2ED2"><w:t xml:space="preserve"> student, a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00AC2ED2"><w:t>Master’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00AC2ED2"><w:t xml:space="preserve"> student (in Data Engineering), and a high-powered sophomore BS student.</w:t></w:r></w:p>'
$p15.InsertXML($x15)

$p19 = $d.Paragraphs(19).Range
$x19 = '<w:p w14:paraId="350DB57E" w14:textId="07FCC551" w:rsidR="007F14BF" w:rsidRDefault="007F14BF" w:rsidP="005C3DB4"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">The AUBRIA Team is passionate about developing AI-Generated Keynote Talks for Conferences, Workshops, Meetings, etc. If you are interested in our team developing an AI-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Geneerated</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Keynote talk for you or your </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>organization</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> contact us!!!</w:t></w:r></w:p>'
$p19.InsertXML($x19)

Write-Output "done"